# Update cryptocurrency price (column D) and 1h volume/change (column E) figures
# for the coin list on the active worksheet, per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '26.180.27'
$ws.Range("E2").Value = '  +0.50%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.658.56'
$ws.Range("E3").Value = '  +0.15%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.37%  '
# Row 5: BNB
$ws.Range("D5").Value = '''215.77'
$ws.Range("E5").Value = '  +4.74%  '
# Row 6: XRP
$ws.Range("D6").Value = '''0.5262'
$ws.Range("E6").Value = '  +2.05%  '
# Row 7: USDC
$ws.Range("E7").Value = '  -0.30%  '
# Row 8: Cardano
$ws.Range("D8").Value = '''0.2633'
$ws.Range("E8").Value = '  +2.24%  '
# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.06406'
$ws.Range("E9").Value = '  +1.97%  '
# Row 10: Solana
$ws.Range("D10").Value = '''20.95'
$ws.Range("E10").Value = '  +0.61%  '
# Row 11: TRON
$ws.Range("D11").Value = '''0.07770'
$ws.Range("E11").Value = '  +3.03%  '
# Row 12: WrappedEther
$ws.Range("D12").Value = '1.655.32'
$ws.Range("E12").Value = '  -0.16%  '
# Row 13: Polkadot
$ws.Range("E13").Value = '  +1.74%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '1.884.57'
$ws.Range("E14").Value = '  -0.12%  '
# Row 15: Polygon
$ws.Range("D15").Value = '''0.5542'
$ws.Range("E15").Value = '  +3.23%  '
# Row 16: ShibaInu
$ws.Range("D16").Value = '0.0₅8291'
$ws.Range("E16").Value = '  +4.84%  '
# Row 17: Litecoin
$ws.Range("D17").Value = '''65.32'
$ws.Range("E17").Value = '  -1.16%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '26.191.97'
$ws.Range("E18").Value = '  +0.49%  '
# Row 19: Dai
$ws.Range("E19").Value = '  -0.30%  '
# Row 20: Uniswap
$ws.Range("D20").Value = '''4.768'
$ws.Range("E20").Value = '  +1.69%  '
# Row 21: BitcoinCash
$ws.Range("D21").Value = '''190.85'
$ws.Range("E21").Value = '  +2.01%  '
# Row 22: Avalanche
$ws.Range("D22").Value = '''10.31'
$ws.Range("E22").Value = '  +1.88%  '
# Row 23: Chainlink
$ws.Range("D23").Value = '''6.372'
$ws.Range("E23").Value = '  +3.39%  '
# Row 24: BinanceUSD
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  -0.44%  '
# Row 25: Monero
$ws.Range("D25").Value = '''143.00'
$ws.Range("E25").Value = '  -3.45%  '
# Row 26: Stellar
$ws.Range("D26").Value = '''0.1261'
$ws.Range("E26").Value = '  +4.12%  '
# Row 27: Cosmos
$ws.Range("D27").Value = '''7.428'
$ws.Range("E27").Value = '  +0.77%  '
# Row 28: EthereumClassic
$ws.Range("D28").Value = '''16.04'
$ws.Range("E28").Value = '  +2.61%  '
# Row 29: Toncoin
$ws.Range("D29").Value = '''1.430'
$ws.Range("E29").Value = '  +2.45%  '
# Row 30: Hedera
$ws.Range("E30").Value = '  +2.55%  '
# Row 31: PancakeSwap
$ws.Range("D31").Value = '''1.267'
$ws.Range("E31").Value = '  +0.95%  '
# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = '''3.558'
$ws.Range("E32").Value = '  +2.76%  '
# Row 33: Filecoin
$ws.Range("D33").Value = '''3.430'
$ws.Range("E33").Value = '  +1.25%  '
# Row 34: LidoDAOToken
$ws.Range("D34").Value = '''1.664'
$ws.Range("E34").Value = '  +1.99%  '
# Row 35: ARBITRUM
$ws.Range("D35").Value = '''1.002'
$ws.Range("E35").Value = '  +2.10%  '
# Row 36: HuobiToken
$ws.Range("E36").Value = '  +0.41%  '
# Row 37: MXToken
$ws.Range("D37").Value = '''2.761'
$ws.Range("E37").Value = '  +0.24%  '
# Row 38: ImmutableX
$ws.Range("D38").Value = '''0.5684'
$ws.Range("E38").Value = '  -3.00%  '
# Row 39: VeChain
$ws.Range("D39").Value = '''0.01607'
$ws.Range("E39").Value = '  +0.95%  '
# Row 40: FraxShare
$ws.Range("D40").Value = '''5.916'
$ws.Range("E40").Value = '  +0.22%  '
# Row 41: TrustWalletToken
$ws.Range("D41").Value = '''0.8552'
$ws.Range("E41").Value = '  +1.33%  '
# Row 42: PaxDollar
$ws.Range("E42").Value = '  -0.25%  '
# Row 43: Maker
$ws.Range("D43").Value = '1.033.50'
$ws.Range("E43").Value = '  -6.07%  '
# Row 44: Quant
$ws.Range("D44").Value = '''99.64'
$ws.Range("E44").Value = '  -0.58%  '
# Row 45: RocketPoolETH
$ws.Range("D45").Value = '1.806.96'
$ws.Range("E45").Value = '  -0.42%  '
# Row 46: BabyDogeCoin
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").Value = '  +1.88%  '
# Row 47: Aave
$ws.Range("D47").Value = '''56.26'
$ws.Range("E47").Value = '  +2.47%  '
# Row 48: Frax
$ws.Range("E48").Value = '  +0.23%  '
# Row 49: EnergySwap
$ws.Range("D49").Value = '''8.069'
$ws.Range("E49").Value = '  +0.86%  '
# Row 50: Cronos
$ws.Range("D50").Value = '''0.05166'
$ws.Range("E50").Value = '  -1.14%  '
# Row 51: Aptos
$ws.Range("D51").Value = '''5.996'
$ws.Range("E51").Value = '  +2.67%  '
